$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.376096248626709
$ws.Range("B1").Value = 2.564220428466797
$ws.Range("C1").Value = 6.469056606292725
$ws.Range("D1").Value = 2.386465311050415
$ws.Range("E1").Value = 1.219548463821411
